# Regenerate the localization-status report for archive:
#  - Status text "Ready for handoff" -> "In Translation" on all sheets
#  - Narrow the "Status" column on each sheet to match the shorter text

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# Update the status values (Overview rolls the per-locale status into E/F,
# the locale sheets carry it in column C).
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# Re-size the status columns to fit the new, shorter text.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
